# Auto-generated edit script: updates market-price columns (H-N)
# on several leve-profit sheets, matching the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2000
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 2000
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H125").Value = 5221.222
$ws.Range("I125").Value = 4287.6
$ws.Range("J125").Value = 5580.3076
$ws.Range("K125").Value = 38588.4
$ws.Range("L125").Value = 50222.7684
$ws.Range("M125").Value = -36128.4
$ws.Range("N125").Value = -55142.7684

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H137").Value = 1866.1538
$ws.Range("I137").Value = 1866.4286
$ws.Range("J137").Value = 1865.8334
$ws.Range("K137").Value = 5599.2858
$ws.Range("L137").Value = 5597.5002
$ws.Range("M137").Value = -3049.2858
$ws.Range("N137").Value = -10697.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 512.3333
$ws.Range("I4").Value = 114.8
$ws.Range("J4").Value = 2500
$ws.Range("K4").Value = 114.8
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 1.200000000000003
$ws.Range("N4").Value = -2732

$ws.Range("H32").Value = 6333.25
$ws.Range("I32").Value = 5266.7
$ws.Range("K32").Value = 5266.7
$ws.Range("M32").Value = -4979.7

$ws.Range("H61").Value = 4354.2046
$ws.Range("I61").Value = 3123.261
$ws.Range("K61").Value = 3123.261
$ws.Range("M61").Value = -2911.261

$ws.Range("H74").Value = 5469.75
$ws.Range("J74").Value = 18209
$ws.Range("L74").Value = 18209
$ws.Range("N74").Value = -19957

$ws.Range("H77").Value = 5469.75
$ws.Range("J77").Value = 18209
$ws.Range("L77").Value = 91045
$ws.Range("N77").Value = -99781

$ws.Range("H122").Value = 2145.4644
$ws.Range("I122").Value = 2104.4211
$ws.Range("K122").Value = 6313.263300000001
$ws.Range("M122").Value = -3863.263300000001

$ws.Range("H132").Value = 5649.7896
$ws.Range("I132").Value = 5290.4
$ws.Range("K132").Value = 15871.2
$ws.Range("M132").Value = -13341.2

$ws.Range("H136").Value = 4354.2046
$ws.Range("I136").Value = 3123.261
$ws.Range("K136").Value = 9369.782999999999
$ws.Range("M136").Value = -6819.782999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1716
$ws.Range("I105").Value = 1508.5714
$ws.Range("K105").Value = 1508.5714
$ws.Range("M105").Value = 238.4286

$ws.Range("H107").Value = 16493.428
$ws.Range("I107").Value = 2742.5
$ws.Range("J107").Value = 98999
$ws.Range("K107").Value = 2742.5
$ws.Range("L107").Value = 98999
$ws.Range("M107").Value = -822.5
$ws.Range("N107").Value = -102839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 918.13336
$ws.Range("I7").Value = 974.7692
$ws.Range("J7").Value = 550
$ws.Range("K7").Value = 974.7692
$ws.Range("L7").Value = 550
$ws.Range("M7").Value = -861.7692
$ws.Range("N7").Value = -776

$ws.Range("H31").Value = 45707
$ws.Range("J31").Value = 4650
$ws.Range("L31").Value = 4650
$ws.Range("N31").Value = -5240

$ws.Range("H34").Value = 45707
$ws.Range("J34").Value = 4650
$ws.Range("L34").Value = 4650
$ws.Range("N34").Value = -5054

$ws.Range("H62").Value = 7754.846
$ws.Range("I62").Value = 6923.6665
$ws.Range("J62").Value = 9625
$ws.Range("K62").Value = 6923.6665
$ws.Range("L62").Value = 9625
$ws.Range("M62").Value = -6299.6665
$ws.Range("N62").Value = -10873

$ws.Range("H65").Value = 7754.846
$ws.Range("I65").Value = 6923.6665
$ws.Range("J65").Value = 9625
$ws.Range("K65").Value = 34618.3325
$ws.Range("L65").Value = 48125
$ws.Range("M65").Value = -31498.3325
$ws.Range("N65").Value = -54365

$ws.Range("H107").Value = 1320.6154
$ws.Range("I107").Value = 1052.2222
$ws.Range("J107").Value = 1924.5
$ws.Range("K107").Value = 1052.2222
$ws.Range("L107").Value = 1924.5
$ws.Range("M107").Value = 867.7778000000001
$ws.Range("N107").Value = -5764.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11475
$ws.Range("I3").Value = 900
$ws.Range("K3").Value = 2700
$ws.Range("M3").Value = -2588

$ws.Range("H56").Value = 12853.571
$ws.Range("I56").Value = 12853.571
$ws.Range("K56").Value = 12853.571
$ws.Range("M56").Value = -12323.571

$ws.Range("H113").Value = 772.2069
$ws.Range("I113").Value = 763
$ws.Range("J113").Value = 777.0526
$ws.Range("K113").Value = 2289
$ws.Range("L113").Value = 2331.1578
$ws.Range("M113").Value = -119
$ws.Range("N113").Value = -6671.1578

$ws.Range("H133").Value = 20962.445
$ws.Range("J133").Value = 14983.25
$ws.Range("L133").Value = 44949.75
$ws.Range("N133").Value = -55069.75

$ws.Range("H134").Value = 5141.35
$ws.Range("I134").Value = 2926.6875
$ws.Range("K134").Value = 8780.0625
$ws.Range("M134").Value = -3710.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 507500
$ws.Range("J29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15580

$ws.Range("H46").Value = 37500
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H102").Value = 35719684
$ws.Range("I102").Value = 2007.28
$ws.Range("J102").Value = 333367000
$ws.Range("K102").Value = 2007.28
$ws.Range("L102").Value = 333367000
$ws.Range("M102").Value = -385.28
$ws.Range("N102").Value = -333370244

$ws.Range("H126").Value = 20044.785
$ws.Range("I126").Value = 22677.25
$ws.Range("K126").Value = 68031.75
$ws.Range("M126").Value = -65561.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 45455156
$ws.Range("I46").Value = 125000500
$ws.Range("J46").Value = 675.8570999999999
$ws.Range("K46").Value = 125000500
$ws.Range("L46").Value = 675.8570999999999
$ws.Range("M46").Value = -125000312
$ws.Range("N46").Value = -1051.8571

$ws.Range("H136").Value = 4666.381
$ws.Range("I136").Value = 4399.6
$ws.Range("K136").Value = 13198.8
$ws.Range("M136").Value = -10648.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3589.5454
$ws.Range("I126").Value = 2997.4
$ws.Range("J126").Value = 4083
$ws.Range("K126").Value = 8992.200000000001
$ws.Range("L126").Value = 12249
$ws.Range("M126").Value = -6522.200000000001
$ws.Range("N126").Value = -17189
